$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: columns C (estoque), D (data_estoque), E (estoque_depois), F (data_movimento) ---
$ws.Cells.Item(4, 3).Value2 = 170
$ws.Cells.Item(4, 5).Value2 = 170
$ws.Cells.Item(4, 4).Value2 = 45884.32056690704
$ws.Cells.Item(4, 6).Value2 = 45883.64247685186
$ws.Cells.Item(5, 3).Value2 = 106
$ws.Cells.Item(5, 5).Value2 = 106
$ws.Cells.Item(5, 4).Value2 = 45884.32056690704
$ws.Cells.Item(5, 6).Value2 = 45883.73616898148
$ws.Cells.Item(22, 3).Value2 = 6
$ws.Cells.Item(22, 5).Value2 = 6
$ws.Cells.Item(22, 4).Value2 = 45884.32052104706
$ws.Cells.Item(22, 6).Value2 = 45883.54162037037
$ws.Cells.Item(23, 3).Value2 = 72
$ws.Cells.Item(23, 5).Value2 = 72
$ws.Cells.Item(23, 4).Value2 = 45884.32052104706
$ws.Cells.Item(23, 6).Value2 = 45883.63241898148
$ws.Cells.Item(24, 6).Value2 = 45883.53421296296
$ws.Cells.Item(33, 3).Value2 = 2645
$ws.Cells.Item(33, 5).Value2 = 2645
$ws.Cells.Item(33, 4).Value2 = 45884.32052104706
$ws.Cells.Item(33, 6).Value2 = 45883.48231481481
$ws.Cells.Item(52, 3).Value2 = 2
$ws.Cells.Item(52, 5).Value2 = 2
$ws.Cells.Item(52, 4).Value2 = 45884.32056690704
$ws.Cells.Item(52, 6).Value2 = 45883.73616898148
$ws.Cells.Item(56, 3).Value2 = 163
$ws.Cells.Item(56, 5).Value2 = 163
$ws.Cells.Item(56, 4).Value2 = 45884.32056690704
$ws.Cells.Item(56, 6).Value2 = 45883.73616898148
$ws.Cells.Item(58, 3).Value2 = 172
$ws.Cells.Item(58, 5).Value2 = 172
$ws.Cells.Item(58, 4).Value2 = 45884.32056690704
$ws.Cells.Item(58, 6).Value2 = 45883.69670138889
$ws.Cells.Item(59, 3).Value2 = 113
$ws.Cells.Item(59, 5).Value2 = 113
$ws.Cells.Item(59, 4).Value2 = 45884.32056690704
$ws.Cells.Item(59, 6).Value2 = 45883.73616898148
$ws.Cells.Item(63, 3).Value2 = 124
$ws.Cells.Item(63, 5).Value2 = 124
$ws.Cells.Item(63, 4).Value2 = 45884.32056690704
$ws.Cells.Item(63, 6).Value2 = 45883.74665509259
$ws.Cells.Item(70, 3).Value2 = 84
$ws.Cells.Item(70, 5).Value2 = 84
$ws.Cells.Item(70, 4).Value2 = 45884.32052104706
$ws.Cells.Item(70, 6).Value2 = 45883.63241898148
$ws.Cells.Item(82, 3).Value2 = 39
$ws.Cells.Item(82, 5).Value2 = 39
$ws.Cells.Item(82, 4).Value2 = 45884.32052104706
$ws.Cells.Item(82, 6).Value2 = 45883.51311342593
$ws.Cells.Item(83, 3).Value2 = 0
$ws.Cells.Item(83, 5).Value2 = 0
$ws.Cells.Item(83, 4).Value2 = 45884.32056690704
$ws.Cells.Item(83, 6).Value2 = 45883.69670138889
$ws.Cells.Item(85, 3).Value2 = 63
$ws.Cells.Item(85, 5).Value2 = 63
$ws.Cells.Item(85, 4).Value2 = 45884.32056690704
$ws.Cells.Item(85, 6).Value2 = 45883.73616898148
$ws.Cells.Item(86, 3).Value2 = 46
$ws.Cells.Item(86, 5).Value2 = 46
$ws.Cells.Item(86, 4).Value2 = 45884.32056690704
$ws.Cells.Item(86, 6).Value2 = 45883.73616898148
$ws.Cells.Item(87, 3).Value2 = 34
$ws.Cells.Item(87, 5).Value2 = 34
$ws.Cells.Item(87, 4).Value2 = 45884.32052104706
$ws.Cells.Item(87, 6).Value2 = 45883.48247685185
$ws.Cells.Item(91, 3).Value2 = 123
$ws.Cells.Item(91, 5).Value2 = 123
$ws.Cells.Item(91, 4).Value2 = 45884.32056690704
$ws.Cells.Item(91, 6).Value2 = 45883.73616898148
$ws.Cells.Item(94, 3).Value2 = 87
$ws.Cells.Item(94, 5).Value2 = 87
$ws.Cells.Item(94, 4).Value2 = 45884.32056690704
$ws.Cells.Item(94, 6).Value2 = 45883.74665509259
$ws.Cells.Item(101, 3).Value2 = 2087
$ws.Cells.Item(101, 5).Value2 = 2087
$ws.Cells.Item(101, 4).Value2 = 45884.32056690704
$ws.Cells.Item(101, 6).Value2 = 45883.75438657407
$ws.Cells.Item(103, 3).Value2 = -5
$ws.Cells.Item(103, 5).Value2 = -5
$ws.Cells.Item(103, 4).Value2 = 45884.32056690704
$ws.Cells.Item(103, 6).Value2 = 45883.73616898148
$ws.Cells.Item(106, 3).Value2 = -4
$ws.Cells.Item(106, 5).Value2 = -4
$ws.Cells.Item(106, 4).Value2 = 45884.32052104706
$ws.Cells.Item(106, 6).Value2 = 45883.48247685185
$ws.Cells.Item(115, 3).Value2 = 63
$ws.Cells.Item(115, 5).Value2 = 63
$ws.Cells.Item(115, 4).Value2 = 45884.32052104706
$ws.Cells.Item(115, 6).Value2 = 45883.63241898148
$ws.Cells.Item(117, 3).Value2 = 673
$ws.Cells.Item(117, 5).Value2 = 673
$ws.Cells.Item(117, 4).Value2 = 45884.32052104706
$ws.Cells.Item(117, 6).Value2 = 45883.51311342593
$ws.Cells.Item(120, 3).Value2 = 23
$ws.Cells.Item(120, 5).Value2 = 23
$ws.Cells.Item(120, 4).Value2 = 45884.32052104706
$ws.Cells.Item(120, 6).Value2 = 45883.53215277778
$ws.Cells.Item(121, 3).Value2 = 512
$ws.Cells.Item(121, 5).Value2 = 512
$ws.Cells.Item(121, 4).Value2 = 45884.32056690704
$ws.Cells.Item(121, 6).Value2 = 45883.73616898148
$ws.Cells.Item(123, 3).Value2 = 110
$ws.Cells.Item(123, 5).Value2 = 110
$ws.Cells.Item(123, 4).Value2 = 45884.32056690704
$ws.Cells.Item(123, 6).Value2 = 45883.73616898148
$ws.Cells.Item(125, 3).Value2 = 735
$ws.Cells.Item(125, 5).Value2 = 735
$ws.Cells.Item(125, 4).Value2 = 45884.32056690704
$ws.Cells.Item(125, 6).Value2 = 45883.74665509259
$ws.Cells.Item(139, 3).Value2 = 97
$ws.Cells.Item(139, 5).Value2 = 97
$ws.Cells.Item(139, 4).Value2 = 45884.32052104706
$ws.Cells.Item(139, 6).Value2 = 45883.48049768519
$ws.Cells.Item(175, 3).Value2 = 122
$ws.Cells.Item(175, 5).Value2 = 122
$ws.Cells.Item(175, 4).Value2 = 45884.32052104706
$ws.Cells.Item(175, 6).Value2 = 45883.51311342593
$ws.Cells.Item(177, 3).Value2 = -1
$ws.Cells.Item(177, 5).Value2 = -1
$ws.Cells.Item(177, 4).Value2 = 45884.32052104706
$ws.Cells.Item(177, 6).Value2 = 45883.54162037037
$ws.Cells.Item(185, 3).Value2 = 90
$ws.Cells.Item(185, 5).Value2 = 90
$ws.Cells.Item(185, 4).Value2 = 45884.32052104706
$ws.Cells.Item(185, 6).Value2 = 45883.48231481481
$ws.Cells.Item(192, 3).Value2 = 51
$ws.Cells.Item(192, 5).Value2 = 51
$ws.Cells.Item(192, 4).Value2 = 45884.32052104706
$ws.Cells.Item(192, 6).Value2 = 45883.63241898148
$ws.Cells.Item(195, 3).Value2 = 178
$ws.Cells.Item(195, 5).Value2 = 178
$ws.Cells.Item(195, 4).Value2 = 45884.32056690704
$ws.Cells.Item(195, 6).Value2 = 45883.64247685186
$ws.Cells.Item(209, 3).Value2 = 24
$ws.Cells.Item(209, 5).Value2 = 24
$ws.Cells.Item(209, 4).Value2 = 45884.32056690704
$ws.Cells.Item(209, 6).Value2 = 45883.73616898148
$ws.Cells.Item(217, 3).Value2 = 43
$ws.Cells.Item(217, 5).Value2 = 43
$ws.Cells.Item(217, 4).Value2 = 45884.32052104706
$ws.Cells.Item(217, 6).Value2 = 45883.63241898148
$ws.Cells.Item(226, 3).Value2 = 51
$ws.Cells.Item(226, 5).Value2 = 51
$ws.Cells.Item(226, 4).Value2 = 45884.32056690704
$ws.Cells.Item(226, 6).Value2 = 45883.66209490741
$ws.Cells.Item(235, 3).Value2 = 111
$ws.Cells.Item(235, 5).Value2 = 111
$ws.Cells.Item(235, 4).Value2 = 45884.32056690704
$ws.Cells.Item(235, 6).Value2 = 45883.73616898148
$ws.Cells.Item(236, 3).Value2 = 149
$ws.Cells.Item(236, 5).Value2 = 149
$ws.Cells.Item(236, 4).Value2 = 45884.32052104706
$ws.Cells.Item(236, 6).Value2 = 45883.48247685185
$ws.Cells.Item(247, 3).Value2 = 209
$ws.Cells.Item(247, 5).Value2 = 209
$ws.Cells.Item(247, 4).Value2 = 45884.32052104706
$ws.Cells.Item(247, 6).Value2 = 45883.48049768519
$ws.Cells.Item(257, 3).Value2 = 40
$ws.Cells.Item(257, 5).Value2 = 40
$ws.Cells.Item(257, 4).Value2 = 45884.32056690704
$ws.Cells.Item(257, 6).Value2 = 45883.66209490741
$ws.Cells.Item(258, 3).Value2 = 47
$ws.Cells.Item(258, 5).Value2 = 47
$ws.Cells.Item(258, 4).Value2 = 45884.32052104706
$ws.Cells.Item(258, 6).Value2 = 45883.63241898148
$ws.Cells.Item(270, 3).Value2 = 161
$ws.Cells.Item(270, 5).Value2 = 161
$ws.Cells.Item(270, 4).Value2 = 45884.32052104706
$ws.Cells.Item(270, 6).Value2 = 45883.42606481481
$ws.Cells.Item(272, 3).Value2 = 344
$ws.Cells.Item(272, 5).Value2 = 344
$ws.Cells.Item(272, 4).Value2 = 45884.32052104706
$ws.Cells.Item(272, 6).Value2 = 45883.54162037037
$ws.Cells.Item(273, 3).Value2 = 97
$ws.Cells.Item(273, 5).Value2 = 97
$ws.Cells.Item(273, 4).Value2 = 45884.32052104706
$ws.Cells.Item(273, 6).Value2 = 45883.38149305555
$ws.Cells.Item(274, 3).Value2 = 23
$ws.Cells.Item(274, 5).Value2 = 23
$ws.Cells.Item(274, 4).Value2 = 45884.32052104706
$ws.Cells.Item(274, 6).Value2 = 45883.48196759259
$ws.Cells.Item(291, 3).Value2 = 44
$ws.Cells.Item(291, 5).Value2 = 44
$ws.Cells.Item(291, 4).Value2 = 45884.32056690704
$ws.Cells.Item(291, 6).Value2 = 45883.69670138889
$ws.Cells.Item(292, 3).Value2 = 138
$ws.Cells.Item(292, 5).Value2 = 138
$ws.Cells.Item(292, 4).Value2 = 45884.32052104706
$ws.Cells.Item(292, 6).Value2 = 45883.48247685185
$ws.Cells.Item(307, 3).Value2 = 5
$ws.Cells.Item(307, 5).Value2 = 5
$ws.Cells.Item(307, 4).Value2 = 45884.32056690704
$ws.Cells.Item(307, 6).Value2 = 45883.73616898148
$ws.Cells.Item(314, 3).Value2 = 593
$ws.Cells.Item(314, 5).Value2 = 593
$ws.Cells.Item(314, 4).Value2 = 45884.32052104706
$ws.Cells.Item(314, 6).Value2 = 45883.42606481481
$ws.Cells.Item(315, 3).Value2 = 184
$ws.Cells.Item(315, 5).Value2 = 184
$ws.Cells.Item(315, 4).Value2 = 45884.32052104706
$ws.Cells.Item(315, 6).Value2 = 45883.48231481481
$ws.Cells.Item(319, 3).Value2 = 107
$ws.Cells.Item(319, 5).Value2 = 107
$ws.Cells.Item(319, 4).Value2 = 45884.32052104706
$ws.Cells.Item(319, 6).Value2 = 45883.51311342593
$ws.Cells.Item(326, 3).Value2 = 404
$ws.Cells.Item(326, 5).Value2 = 404
$ws.Cells.Item(326, 4).Value2 = 45884.32056690704
$ws.Cells.Item(326, 6).Value2 = 45883.69237268518
$ws.Cells.Item(328, 3).Value2 = 154
$ws.Cells.Item(328, 5).Value2 = 154
$ws.Cells.Item(328, 4).Value2 = 45884.32052104706
$ws.Cells.Item(328, 6).Value2 = 45883.63241898148
$ws.Cells.Item(334, 3).Value2 = 75
$ws.Cells.Item(334, 5).Value2 = 75
$ws.Cells.Item(334, 4).Value2 = 45884.32056690704
$ws.Cells.Item(334, 6).Value2 = 45883.64247685186
$ws.Cells.Item(354, 3).Value2 = 152
$ws.Cells.Item(354, 5).Value2 = 152
$ws.Cells.Item(354, 4).Value2 = 45884.32052104706
$ws.Cells.Item(354, 6).Value2 = 45883.63241898148
$ws.Cells.Item(358, 3).Value2 = 63
$ws.Cells.Item(358, 5).Value2 = 63
$ws.Cells.Item(358, 4).Value2 = 45884.32056690704
$ws.Cells.Item(358, 6).Value2 = 45883.645625
$ws.Cells.Item(363, 3).Value2 = -2
$ws.Cells.Item(363, 5).Value2 = -2
$ws.Cells.Item(363, 4).Value2 = 45884.32056690704
$ws.Cells.Item(363, 6).Value2 = 45883.66209490741
$ws.Cells.Item(387, 3).Value2 = 51
$ws.Cells.Item(387, 5).Value2 = 51
$ws.Cells.Item(387, 4).Value2 = 45884.32056690704
$ws.Cells.Item(387, 6).Value2 = 45883.64517361111
$ws.Cells.Item(390, 3).Value2 = 96
$ws.Cells.Item(390, 5).Value2 = 96
$ws.Cells.Item(390, 4).Value2 = 45884.32056690704
$ws.Cells.Item(390, 6).Value2 = 45883.74665509259
$ws.Cells.Item(395, 3).Value2 = 2
$ws.Cells.Item(395, 5).Value2 = 2
$ws.Cells.Item(395, 4).Value2 = 45884.32056690704
$ws.Cells.Item(395, 6).Value2 = 45883.73616898148
$ws.Cells.Item(404, 3).Value2 = 59
$ws.Cells.Item(404, 5).Value2 = 59
$ws.Cells.Item(404, 4).Value2 = 45884.32056690704
$ws.Cells.Item(404, 6).Value2 = 45883.645625
$ws.Cells.Item(408, 3).Value2 = 122
$ws.Cells.Item(408, 5).Value2 = 122
$ws.Cells.Item(408, 4).Value2 = 45884.32056690704
$ws.Cells.Item(408, 6).Value2 = 45883.645625
$ws.Cells.Item(410, 3).Value2 = 2165
$ws.Cells.Item(410, 5).Value2 = 2165
$ws.Cells.Item(410, 4).Value2 = 45884.32052104706
$ws.Cells.Item(410, 6).Value2 = 45883.48247685185
$ws.Cells.Item(414, 3).Value2 = 36
$ws.Cells.Item(414, 5).Value2 = 36
$ws.Cells.Item(414, 4).Value2 = 45884.32056690704
$ws.Cells.Item(414, 6).Value2 = 45883.74665509259
$ws.Cells.Item(416, 3).Value2 = 12
$ws.Cells.Item(416, 5).Value2 = 12
$ws.Cells.Item(416, 4).Value2 = 45884.32052104706
$ws.Cells.Item(416, 6).Value2 = 45883.48049768519
$ws.Cells.Item(418, 3).Value2 = 442
$ws.Cells.Item(418, 5).Value2 = 442
$ws.Cells.Item(418, 4).Value2 = 45884.32052104706
$ws.Cells.Item(418, 6).Value2 = 45883.48247685185
$ws.Cells.Item(425, 3).Value2 = 516
$ws.Cells.Item(425, 5).Value2 = 516
$ws.Cells.Item(425, 4).Value2 = 45884.32052104706
$ws.Cells.Item(425, 6).Value2 = 45883.53215277778
$ws.Cells.Item(435, 3).Value2 = 75
$ws.Cells.Item(435, 5).Value2 = 75
$ws.Cells.Item(435, 4).Value2 = 45884.32052104706
$ws.Cells.Item(435, 6).Value2 = 45883.48023148148
$ws.Cells.Item(437, 3).Value2 = 40
$ws.Cells.Item(437, 5).Value2 = 40
$ws.Cells.Item(437, 4).Value2 = 45884.32052104706
$ws.Cells.Item(437, 6).Value2 = 45883.63241898148
$ws.Cells.Item(439, 3).Value2 = 202
$ws.Cells.Item(439, 5).Value2 = 202
$ws.Cells.Item(439, 4).Value2 = 45884.32052104706
$ws.Cells.Item(439, 6).Value2 = 45883.63241898148
$ws.Cells.Item(443, 3).Value2 = 7
$ws.Cells.Item(443, 5).Value2 = 7
$ws.Cells.Item(443, 4).Value2 = 45884.32052104706
$ws.Cells.Item(443, 6).Value2 = 45883.47914351852
$ws.Cells.Item(461, 3).Value2 = 15
$ws.Cells.Item(461, 5).Value2 = 15
$ws.Cells.Item(461, 4).Value2 = 45884.32052104706
$ws.Cells.Item(461, 6).Value2 = 45883.48247685185
$ws.Cells.Item(469, 3).Value2 = 2456
$ws.Cells.Item(469, 5).Value2 = 2456
$ws.Cells.Item(469, 4).Value2 = 45884.32056690704
$ws.Cells.Item(469, 6).Value2 = 45883.73616898148
$ws.Cells.Item(472, 3).Value2 = 76
$ws.Cells.Item(472, 5).Value2 = 76
$ws.Cells.Item(472, 4).Value2 = 45884.32056690704
$ws.Cells.Item(472, 6).Value2 = 45883.66209490741
$ws.Cells.Item(480, 3).Value2 = 162
$ws.Cells.Item(480, 5).Value2 = 162
$ws.Cells.Item(480, 4).Value2 = 45884.32056690704
$ws.Cells.Item(480, 6).Value2 = 45883.70675925926
$ws.Cells.Item(510, 3).Value2 = 20
$ws.Cells.Item(510, 5).Value2 = 20
$ws.Cells.Item(510, 4).Value2 = 45884.32052104706
$ws.Cells.Item(510, 6).Value2 = 45883.38149305555
$ws.Cells.Item(525, 3).Value2 = 14
$ws.Cells.Item(525, 5).Value2 = 14
$ws.Cells.Item(525, 4).Value2 = 45884.32052104706
$ws.Cells.Item(525, 6).Value2 = 45883.48980324074
$ws.Cells.Item(528, 3).Value2 = 240
$ws.Cells.Item(528, 5).Value2 = 240
$ws.Cells.Item(528, 4).Value2 = 45884.32056690704
$ws.Cells.Item(528, 6).Value2 = 45883.645625
$ws.Cells.Item(535, 3).Value2 = 388
$ws.Cells.Item(535, 5).Value2 = 388
$ws.Cells.Item(535, 4).Value2 = 45884.32052104706
$ws.Cells.Item(535, 6).Value2 = 45883.48049768519
$ws.Cells.Item(545, 3).Value2 = 23
$ws.Cells.Item(545, 5).Value2 = 23
$ws.Cells.Item(545, 4).Value2 = 45884.32052104706
$ws.Cells.Item(545, 6).Value2 = 45883.48247685185
$ws.Cells.Item(569, 3).Value2 = 126
$ws.Cells.Item(569, 5).Value2 = 126
$ws.Cells.Item(569, 4).Value2 = 45884.32056690704
$ws.Cells.Item(569, 6).Value2 = 45883.66209490741
$ws.Cells.Item(581, 3).Value2 = 12
$ws.Cells.Item(581, 5).Value2 = 12
$ws.Cells.Item(581, 4).Value2 = 45884.32052104706
$ws.Cells.Item(581, 6).Value2 = 45883.53533564815
$ws.Cells.Item(613, 3).Value2 = -8
$ws.Cells.Item(613, 5).Value2 = -8
$ws.Cells.Item(613, 4).Value2 = 45884.32052104706
$ws.Cells.Item(613, 6).Value2 = 45883.63241898148
$ws.Cells.Item(634, 3).Value2 = -9
$ws.Cells.Item(634, 5).Value2 = -9
$ws.Cells.Item(634, 4).Value2 = 45884.32052104706
$ws.Cells.Item(634, 6).Value2 = 45883.48247685185
$ws.Cells.Item(652, 3).Value2 = 4
$ws.Cells.Item(652, 5).Value2 = 4
$ws.Cells.Item(652, 4).Value2 = 45884.32052104706
$ws.Cells.Item(652, 6).Value2 = 45883.48049768519
$ws.Cells.Item(657, 3).Value2 = 933
$ws.Cells.Item(657, 5).Value2 = 933
$ws.Cells.Item(657, 4).Value2 = 45884.32056690704
$ws.Cells.Item(657, 6).Value2 = 45883.645625
$ws.Cells.Item(660, 3).Value2 = 276
$ws.Cells.Item(660, 5).Value2 = 276
$ws.Cells.Item(660, 4).Value2 = 45884.32052104706
$ws.Cells.Item(660, 6).Value2 = 45883.4808912037
$ws.Cells.Item(679, 4).Value2 = 45884.32056690704
$ws.Cells.Item(681, 3).Value2 = 102
$ws.Cells.Item(681, 5).Value2 = 102
$ws.Cells.Item(681, 4).Value2 = 45884.32052104706
$ws.Cells.Item(681, 6).Value2 = 45883.47935185185
$ws.Cells.Item(689, 3).Value2 = 43
$ws.Cells.Item(689, 5).Value2 = 43
$ws.Cells.Item(689, 4).Value2 = 45884.32056690704
$ws.Cells.Item(689, 6).Value2 = 45883.6790625
$ws.Cells.Item(701, 3).Value2 = 8
$ws.Cells.Item(701, 5).Value2 = 8
$ws.Cells.Item(701, 4).Value2 = 45884.32052104706
$ws.Cells.Item(701, 6).Value2 = 45883.4808912037
$ws.Cells.Item(710, 3).Value2 = 45
$ws.Cells.Item(710, 5).Value2 = 45
$ws.Cells.Item(710, 4).Value2 = 45884.32052104706
$ws.Cells.Item(710, 6).Value2 = 45883.4807175926
$ws.Cells.Item(720, 3).Value2 = 254
$ws.Cells.Item(720, 5).Value2 = 254
$ws.Cells.Item(720, 4).Value2 = 45884.32056690704
$ws.Cells.Item(720, 6).Value2 = 45883.73616898148
$ws.Cells.Item(721, 3).Value2 = 9
$ws.Cells.Item(721, 5).Value2 = 9
$ws.Cells.Item(721, 4).Value2 = 45884.32052104706
$ws.Cells.Item(721, 6).Value2 = 45883.51311342593
$ws.Cells.Item(729, 3).Value2 = 6
$ws.Cells.Item(729, 5).Value2 = 6
$ws.Cells.Item(729, 4).Value2 = 45884.32052104706
$ws.Cells.Item(729, 6).Value2 = 45883.4808912037
$ws.Cells.Item(730, 3).Value2 = 32
$ws.Cells.Item(730, 5).Value2 = 32
$ws.Cells.Item(730, 4).Value2 = 45884.32052104706
$ws.Cells.Item(730, 6).Value2 = 45883.4808912037
$ws.Cells.Item(731, 3).Value2 = 0
$ws.Cells.Item(731, 5).Value2 = 0
$ws.Cells.Item(731, 4).Value2 = 45884.32052104706
$ws.Cells.Item(731, 6).Value2 = 45883.48231481481
$ws.Cells.Item(732, 3).Value2 = 976
$ws.Cells.Item(732, 5).Value2 = 976
$ws.Cells.Item(732, 4).Value2 = 45884.32052104706
$ws.Cells.Item(732, 6).Value2 = 45883.48196759259
$ws.Cells.Item(735, 3).Value2 = 11
$ws.Cells.Item(735, 5).Value2 = 11
$ws.Cells.Item(735, 4).Value2 = 45884.32052104706
$ws.Cells.Item(735, 6).Value2 = 45883.54490740741
$ws.Cells.Item(737, 3).Value2 = 38
$ws.Cells.Item(737, 5).Value2 = 38
$ws.Cells.Item(737, 4).Value2 = 45884.32052104706
$ws.Cells.Item(737, 6).Value2 = 45883.54162037037
$ws.Cells.Item(741, 3).Value2 = 36
$ws.Cells.Item(741, 5).Value2 = 36
$ws.Cells.Item(741, 4).Value2 = 45884.32052104706
$ws.Cells.Item(741, 6).Value2 = 45883.48247685185
$ws.Cells.Item(743, 3).Value2 = 36
$ws.Cells.Item(743, 5).Value2 = 36
$ws.Cells.Item(743, 4).Value2 = 45884.32056690704
$ws.Cells.Item(743, 6).Value2 = 45883.66178240741
$ws.Cells.Item(748, 4).Value2 = 45884.32056690704
$ws.Cells.Item(754, 3).Value2 = -3
$ws.Cells.Item(754, 5).Value2 = -3
$ws.Cells.Item(754, 4).Value2 = 45884.32052104706
$ws.Cells.Item(754, 6).Value2 = 45883.54162037037
$ws.Cells.Item(757, 3).Value2 = 78
$ws.Cells.Item(757, 5).Value2 = 78
$ws.Cells.Item(757, 4).Value2 = 45884.32052104706
$ws.Cells.Item(757, 6).Value2 = 45883.63241898148
$ws.Cells.Item(763, 3).Value2 = 2
$ws.Cells.Item(763, 5).Value2 = 2
$ws.Cells.Item(763, 4).Value2 = 45884.32052104706
$ws.Cells.Item(763, 6).Value2 = 45883.42606481481
$ws.Cells.Item(771, 3).Value2 = 12
$ws.Cells.Item(771, 5).Value2 = 12
$ws.Cells.Item(771, 4).Value2 = 45884.32052104706
$ws.Cells.Item(771, 6).Value2 = 45883.42606481481
$ws.Cells.Item(774, 3).Value2 = 97
$ws.Cells.Item(774, 5).Value2 = 97
$ws.Cells.Item(774, 4).Value2 = 45884.32056690704
$ws.Cells.Item(774, 6).Value2 = 45883.66209490741
$ws.Cells.Item(776, 3).Value2 = 1616
$ws.Cells.Item(776, 5).Value2 = 1616
$ws.Cells.Item(776, 4).Value2 = 45884.32052104706
$ws.Cells.Item(776, 6).Value2 = 45883.51311342593
$ws.Cells.Item(778, 3).Value2 = 86
$ws.Cells.Item(778, 5).Value2 = 86
$ws.Cells.Item(778, 4).Value2 = 45884.32052104706
$ws.Cells.Item(778, 6).Value2 = 45883.51311342593
$ws.Cells.Item(807, 3).Value2 = -3
$ws.Cells.Item(807, 5).Value2 = -3
$ws.Cells.Item(807, 4).Value2 = 45884.32052104706
$ws.Cells.Item(807, 6).Value2 = 45883.54162037037
$ws.Cells.Item(810, 3).Value2 = 1
$ws.Cells.Item(810, 5).Value2 = 1
$ws.Cells.Item(810, 4).Value2 = 45884.32052104706
$ws.Cells.Item(810, 6).Value2 = 45883.5344212963
$ws.Cells.Item(812, 3).Value2 = 21
$ws.Cells.Item(812, 5).Value2 = 21
$ws.Cells.Item(812, 4).Value2 = 45884.32052104706
$ws.Cells.Item(812, 6).Value2 = 45883.47771990741
$ws.Cells.Item(816, 3).Value2 = 8
$ws.Cells.Item(816, 5).Value2 = 8
$ws.Cells.Item(816, 4).Value2 = 45884.32052104706
$ws.Cells.Item(816, 6).Value2 = 45883.4808912037
$ws.Cells.Item(824, 3).Value2 = 45
$ws.Cells.Item(824, 5).Value2 = 45
$ws.Cells.Item(824, 4).Value2 = 45884.32056690704
$ws.Cells.Item(824, 6).Value2 = 45883.73616898148
$ws.Cells.Item(826, 3).Value2 = 111
$ws.Cells.Item(826, 5).Value2 = 111
$ws.Cells.Item(826, 4).Value2 = 45884.32056690704
$ws.Cells.Item(826, 6).Value2 = 45883.74665509259
$ws.Cells.Item(839, 3).Value2 = 5
$ws.Cells.Item(839, 5).Value2 = 5
$ws.Cells.Item(839, 4).Value2 = 45884.32056690704
$ws.Cells.Item(839, 6).Value2 = 45883.73616898148
$ws.Cells.Item(854, 3).Value2 = 132
$ws.Cells.Item(854, 5).Value2 = 132
$ws.Cells.Item(854, 4).Value2 = 45884.32052104706
$ws.Cells.Item(854, 6).Value2 = 45883.38149305555
$ws.Cells.Item(871, 3).Value2 = 31
$ws.Cells.Item(871, 5).Value2 = 31
$ws.Cells.Item(871, 4).Value2 = 45884.32056690704
$ws.Cells.Item(871, 6).Value2 = 45883.73616898148
$ws.Cells.Item(872, 3).Value2 = 378
$ws.Cells.Item(872, 5).Value2 = 378
$ws.Cells.Item(872, 4).Value2 = 45884.32056690704
$ws.Cells.Item(872, 6).Value2 = 45883.645625
$ws.Cells.Item(876, 3).Value2 = 780
$ws.Cells.Item(876, 5).Value2 = 780
$ws.Cells.Item(876, 4).Value2 = 45884.32056690704
$ws.Cells.Item(876, 6).Value2 = 45883.645625
$ws.Cells.Item(883, 3).Value2 = 370
$ws.Cells.Item(883, 5).Value2 = 370
$ws.Cells.Item(883, 4).Value2 = 45884.32052104706
$ws.Cells.Item(883, 6).Value2 = 45883.48247685185
$ws.Cells.Item(889, 3).Value2 = 35
$ws.Cells.Item(889, 5).Value2 = 35
$ws.Cells.Item(889, 4).Value2 = 45884.32052104706
$ws.Cells.Item(889, 6).Value2 = 45883.63241898148
$ws.Cells.Item(924, 3).Value2 = -50
$ws.Cells.Item(924, 5).Value2 = -50
$ws.Cells.Item(924, 4).Value2 = 45884.32056690704
$ws.Cells.Item(924, 6).Value2 = 45883.74665509259
$ws.Cells.Item(941, 3).Value2 = 1
$ws.Cells.Item(941, 5).Value2 = 1
$ws.Cells.Item(941, 4).Value2 = 45884.32056690704
$ws.Cells.Item(941, 6).Value2 = 45883.63241898148
$ws.Cells.Item(962, 3).Value2 = 582
$ws.Cells.Item(962, 5).Value2 = 582
$ws.Cells.Item(962, 4).Value2 = 45884.32056690704
$ws.Cells.Item(962, 6).Value2 = 45883.63241898148
$ws.Cells.Item(963, 3).Value2 = 857
$ws.Cells.Item(963, 5).Value2 = 857
$ws.Cells.Item(963, 4).Value2 = 45884.32056690704
$ws.Cells.Item(963, 6).Value2 = 45883.63241898148
$ws.Cells.Item(977, 3).Value2 = 73
$ws.Cells.Item(977, 5).Value2 = 73
$ws.Cells.Item(977, 4).Value2 = 45884.32052104706
$ws.Cells.Item(977, 6).Value2 = 45883.48049768519
$ws.Cells.Item(1025, 3).Value2 = 130
$ws.Cells.Item(1025, 5).Value2 = 130
$ws.Cells.Item(1025, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1025, 6).Value2 = 45883.73616898148
$ws.Cells.Item(1027, 3).Value2 = 8
$ws.Cells.Item(1027, 5).Value2 = 8
$ws.Cells.Item(1027, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1027, 6).Value2 = 45883.63241898148
$ws.Cells.Item(1039, 3).Value2 = 433
$ws.Cells.Item(1039, 5).Value2 = 433
$ws.Cells.Item(1039, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1039, 6).Value2 = 45883.645625
$ws.Cells.Item(1062, 3).Value2 = 280
$ws.Cells.Item(1062, 5).Value2 = 280
$ws.Cells.Item(1062, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1062, 6).Value2 = 45883.64247685186
$ws.Cells.Item(1065, 3).Value2 = 2
$ws.Cells.Item(1065, 5).Value2 = 2
$ws.Cells.Item(1065, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1065, 6).Value2 = 45883.63241898148
$ws.Cells.Item(1094, 3).Value2 = 13
$ws.Cells.Item(1094, 5).Value2 = 13
$ws.Cells.Item(1094, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1094, 6).Value2 = 45883.74665509259
$ws.Cells.Item(1113, 3).Value2 = -4
$ws.Cells.Item(1113, 5).Value2 = -4
$ws.Cells.Item(1113, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1113, 6).Value2 = 45883.51311342593
$ws.Cells.Item(1126, 3).Value2 = 367
$ws.Cells.Item(1126, 5).Value2 = 367
$ws.Cells.Item(1126, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1126, 6).Value2 = 45883.51311342593
$ws.Cells.Item(1149, 3).Value2 = 20
$ws.Cells.Item(1149, 5).Value2 = 20
$ws.Cells.Item(1149, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1149, 6).Value2 = 45883.38149305555
$ws.Cells.Item(1191, 3).Value2 = 0
$ws.Cells.Item(1191, 5).Value2 = 0
$ws.Cells.Item(1191, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1191, 6).Value2 = 45883.46005787037
$ws.Cells.Item(1193, 3).Value2 = 87
$ws.Cells.Item(1193, 5).Value2 = 87
$ws.Cells.Item(1193, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1193, 6).Value2 = 45883.51311342593
$ws.Cells.Item(1199, 3).Value2 = -2
$ws.Cells.Item(1199, 5).Value2 = -2
$ws.Cells.Item(1199, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1199, 6).Value2 = 45883.69237268518
$ws.Cells.Item(1223, 3).Value2 = 89
$ws.Cells.Item(1223, 5).Value2 = 89
$ws.Cells.Item(1223, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1223, 6).Value2 = 45883.63241898148
$ws.Cells.Item(1253, 3).Value2 = 238
$ws.Cells.Item(1253, 5).Value2 = 238
$ws.Cells.Item(1253, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1253, 6).Value2 = 45883.70675925926
$ws.Cells.Item(1261, 3).Value2 = -79
$ws.Cells.Item(1261, 5).Value2 = -79
$ws.Cells.Item(1261, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1261, 6).Value2 = 45883.74665509259
$ws.Cells.Item(1266, 3).Value2 = 0
$ws.Cells.Item(1266, 5).Value2 = 0
$ws.Cells.Item(1266, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1266, 6).Value2 = 45883.44671296296
$ws.Cells.Item(1301, 3).Value2 = -5
$ws.Cells.Item(1301, 5).Value2 = -5
$ws.Cells.Item(1301, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1301, 6).Value2 = 45883.54162037037
$ws.Cells.Item(1330, 3).Value2 = 80
$ws.Cells.Item(1330, 5).Value2 = 80
$ws.Cells.Item(1330, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1330, 6).Value2 = 45883.51311342593
$ws.Cells.Item(1342, 3).Value2 = 745
$ws.Cells.Item(1342, 5).Value2 = 745
$ws.Cells.Item(1342, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1342, 6).Value2 = 45883.51311342593
$ws.Cells.Item(1385, 3).Value2 = 115
$ws.Cells.Item(1385, 5).Value2 = 115
$ws.Cells.Item(1385, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1385, 6).Value2 = 45883.38149305555
$ws.Cells.Item(1390, 3).Value2 = 18
$ws.Cells.Item(1390, 5).Value2 = 18
$ws.Cells.Item(1390, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1390, 6).Value2 = 45883.74665509259
$ws.Cells.Item(1394, 6).Value2 = 45883.63748842593
$ws.Cells.Item(1403, 3).Value2 = 64
$ws.Cells.Item(1403, 5).Value2 = 64
$ws.Cells.Item(1403, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1403, 6).Value2 = 45883.53215277778
$ws.Cells.Item(1417, 3).Value2 = 1
$ws.Cells.Item(1417, 5).Value2 = 1
$ws.Cells.Item(1417, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1417, 6).Value2 = 45883.42606481481
$ws.Cells.Item(1423, 3).Value2 = 28
$ws.Cells.Item(1423, 5).Value2 = 28
$ws.Cells.Item(1423, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1423, 6).Value2 = 45883.73616898148
$ws.Cells.Item(1428, 3).Value2 = 141
$ws.Cells.Item(1428, 5).Value2 = 141
$ws.Cells.Item(1428, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1428, 6).Value2 = 45883.48196759259
$ws.Cells.Item(1446, 3).Value2 = 41
$ws.Cells.Item(1446, 5).Value2 = 41
$ws.Cells.Item(1446, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1446, 6).Value2 = 45883.38149305555
$ws.Cells.Item(1507, 3).Value2 = 19
$ws.Cells.Item(1507, 5).Value2 = 19
$ws.Cells.Item(1507, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1507, 6).Value2 = 45883.48247685185
$ws.Cells.Item(1515, 3).Value2 = 0
$ws.Cells.Item(1515, 5).Value2 = 0
$ws.Cells.Item(1515, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1515, 6).Value2 = 45883.53210648148
$ws.Cells.Item(1594, 3).Value2 = 81
$ws.Cells.Item(1594, 5).Value2 = 81
$ws.Cells.Item(1594, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1594, 6).Value2 = 45883.74665509259
$ws.Cells.Item(1597, 3).Value2 = 9628
$ws.Cells.Item(1597, 5).Value2 = 9628
$ws.Cells.Item(1597, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1597, 6).Value2 = 45883.63241898148
$ws.Cells.Item(1625, 3).Value2 = 9
$ws.Cells.Item(1625, 5).Value2 = 9
$ws.Cells.Item(1625, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1625, 6).Value2 = 45883.64247685186
$ws.Cells.Item(1635, 3).Value2 = 182
$ws.Cells.Item(1635, 5).Value2 = 182
$ws.Cells.Item(1635, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1635, 6).Value2 = 45883.53215277778
$ws.Cells.Item(1658, 3).Value2 = 2
$ws.Cells.Item(1658, 5).Value2 = 2
$ws.Cells.Item(1658, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1658, 6).Value2 = 45883.53421296296
$ws.Cells.Item(1796, 3).Value2 = 6
$ws.Cells.Item(1796, 5).Value2 = 6
$ws.Cells.Item(1796, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1796, 6).Value2 = 45883.73616898148
$ws.Cells.Item(1810, 3).Value2 = -2
$ws.Cells.Item(1810, 5).Value2 = -2
$ws.Cells.Item(1810, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1810, 6).Value2 = 45883.63241898148
$ws.Cells.Item(1909, 6).Value2 = 45883.65461805555
$ws.Cells.Item(1911, 3).Value2 = 0
$ws.Cells.Item(1911, 5).Value2 = 0
$ws.Cells.Item(1911, 4).Value2 = 45884.32052104706
$ws.Cells.Item(1911, 6).Value2 = 45883.53398148148
$ws.Cells.Item(1932, 3).Value2 = 12
$ws.Cells.Item(1932, 5).Value2 = 12
$ws.Cells.Item(1932, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1932, 6).Value2 = 45883.73616898148
$ws.Cells.Item(1991, 3).Value2 = 1
$ws.Cells.Item(1991, 5).Value2 = 1
$ws.Cells.Item(1991, 4).Value2 = 45884.32056690704
$ws.Cells.Item(1991, 6).Value2 = 45883.73616898148
$ws.Cells.Item(2004, 3).Value2 = 28
$ws.Cells.Item(2004, 5).Value2 = 28
$ws.Cells.Item(2004, 4).Value2 = 45884.32056690704
$ws.Cells.Item(2004, 6).Value2 = 45883.73616898148
$ws.Cells.Item(2024, 3).Value2 = 21
$ws.Cells.Item(2024, 5).Value2 = 21
$ws.Cells.Item(2024, 4).Value2 = 45884.32056690704
$ws.Cells.Item(2024, 6).Value2 = 45883.69670138889
$ws.Cells.Item(2033, 3).Value2 = 11
$ws.Cells.Item(2033, 5).Value2 = 11
$ws.Cells.Item(2033, 4).Value2 = 45884.32052104706
$ws.Cells.Item(2033, 6).Value2 = 45883.44905092593
$ws.Cells.Item(2035, 3).Value2 = 14
$ws.Cells.Item(2035, 5).Value2 = 14
$ws.Cells.Item(2035, 4).Value2 = 45884.32052104706
$ws.Cells.Item(2035, 6).Value2 = 45883.54253472222
$ws.Cells.Item(2042, 3).Value2 = 1
$ws.Cells.Item(2042, 5).Value2 = 1
$ws.Cells.Item(2042, 4).Value2 = 45884.32056690704
$ws.Cells.Item(2042, 6).Value2 = 45883.63241898148
$ws.Cells.Item(2065, 4).Value2 = 45884.32056690704
$ws.Cells.Item(2180, 3).Value2 = 8
$ws.Cells.Item(2180, 5).Value2 = 8
$ws.Cells.Item(2180, 4).Value2 = 45884.32056690704
$ws.Cells.Item(2180, 6).Value2 = 45883.645625
$ws.Cells.Item(2349, 3).Value2 = 0
$ws.Cells.Item(2349, 5).Value2 = 0
$ws.Cells.Item(2349, 4).Value2 = 45884.32052104706
$ws.Cells.Item(2349, 6).Value2 = 45883.53322916666
$ws.Cells.Item(2481, 3).Value2 = 604
$ws.Cells.Item(2481, 5).Value2 = 604
$ws.Cells.Item(2481, 4).Value2 = 45884.32056690704
$ws.Cells.Item(2481, 6).Value2 = 45883.66209490741
$ws.Cells.Item(2484, 3).Value2 = 2118
$ws.Cells.Item(2484, 5).Value2 = 2118
$ws.Cells.Item(2484, 4).Value2 = 45884.32056690704
$ws.Cells.Item(2484, 6).Value2 = 45883.66209490741
$ws.Cells.Item(2602, 3).Value2 = 2
$ws.Cells.Item(2602, 5).Value2 = 2
$ws.Cells.Item(2602, 4).Value2 = 45884.32052104706
$ws.Cells.Item(2602, 6).Value2 = 45883.53931712963
$ws.Cells.Item(2603, 6).Value2 = 45883.54123842593
$ws.Cells.Item(2636, 3).Value2 = 1
$ws.Cells.Item(2636, 5).Value2 = 1
$ws.Cells.Item(2636, 4).Value2 = 45884.32052104706
$ws.Cells.Item(2636, 6).Value2 = 45883.59208333334
$ws.Cells.Item(2641, 3).Value2 = 0
$ws.Cells.Item(2641, 5).Value2 = 0
$ws.Cells.Item(2641, 4).Value2 = 45884.32052104706
$ws.Cells.Item(2641, 6).Value2 = 45883.45003472222
# --- Append new rows 2687 and 2688 ---
$ws.Cells.Item(2687, 1).Value2 = 44807146
$ws.Cells.Item(2687, 2).Value2 = 1
$ws.Cells.Item(2687, 3).Value2 = -30
$ws.Cells.Item(2687, 4).Value2 = 45884.32056690704
$ws.Cells.Item(2687, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2687, 5).Value2 = -30
$ws.Cells.Item(2687, 6).Value2 = 45883.63748842593
$ws.Cells.Item(2687, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2687, 7).Value2 = 0
$ws.Cells.Item(2687, 8).Value2 = "Consistente"

$ws.Cells.Item(2688, 1).Value2 = 44809973
$ws.Cells.Item(2688, 2).Value2 = 1
$ws.Cells.Item(2688, 3).Value2 = 0
$ws.Cells.Item(2688, 4).Value2 = 45884.32056690704
$ws.Cells.Item(2688, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2688, 5).Value2 = 0
$ws.Cells.Item(2688, 7).Value2 = 0
$ws.Cells.Item(2688, 8).Value2 = "Consistente"
